$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '66.015.61'
Set-TextValue 'E2' '  -0.82%  '
Set-TextValue 'D3' '3.301.36'
Set-TextValue 'E3' '  -1.64%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '188.08'
Set-TextValue 'E5' '  +2.66%  '
Set-TextValue 'D6' '554.84'
Set-TextValue 'E6' '  -0.31%  '
Set-TextValue 'D8' '0.585'
Set-TextValue 'E8' '  -2.24%  '
Set-TextValue 'D9' '3.293.51'
Set-TextValue 'E9' '  -1.66%  '
Set-TextValue 'D10' '0.184'
Set-TextValue 'E10' '  -1.74%  '
Set-TextValue 'E11' '  -1.15%  '
Set-TextValue 'D12' '47.43'
Set-TextValue 'E12' '  -0.74%  '
Set-TextValue 'D13' '0.0000270'
Set-TextValue 'E13' '  +0.99%  '
Set-TextValue 'D14' '8.63'
Set-TextValue 'E14' '  -0.74%  '
Set-TextValue 'D15' '3.833.86'
Set-TextValue 'E15' '  -1.33%  '
Set-TextValue 'D16' '614.28'
Set-TextValue 'E16' '  +2.22%  '
Set-TextValue 'D17' '18.04'
Set-TextValue 'E17' '  +1.07%  '
Set-TextValue 'D18' '65.977.73'
Set-TextValue 'E18' '  -0.62%  '
Set-TextValue 'E19' '  -0.13%  '
Set-TextValue 'D20' '3.304.65'
Set-TextValue 'E20' '  -0.89%  '
Set-TextValue 'D21' '10.94'
Set-TextValue 'E21' '  -6.13%  '
Set-TextValue 'D22' '0.907'
Set-TextValue 'E22' '  -0.24%  '
Set-TextValue 'D23' '18.35'
Set-TextValue 'E23' '  +9.15%  '
Set-TextValue 'D24' '101.86'
Set-TextValue 'E24' '  +3.99%  '
Set-TextValue 'E25' '  -1.81%  '
Set-TextValue 'D26' '3.93'
Set-TextValue 'E26' '  -3.14%  '
Set-TextValue 'D28' '2.73'
Set-TextValue 'E28' '  -0.46%  '
Set-TextValue 'D29' '9.60'
Set-TextValue 'E29' '  +1.84%  '
Set-TextValue 'D30' '8.63'
Set-TextValue 'E30' '  -1.91%  '
Set-TextValue 'D31' '30.20'
Set-TextValue 'E31' '  -1.58%  '
Set-TextValue 'D32' '4.06'
Set-TextValue 'E32' '  +5.64%  '
Set-TextValue 'D33' '6.49'
Set-TextValue 'E33' '  +2.54%  '
Set-TextValue 'D34' '559.27'
Set-TextValue 'E34' '  +5.30%  '
Set-TextValue 'D35' '11.05'
Set-TextValue 'E35' '  -0.86%  '
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'D37' '3.824.26'
Set-TextValue 'E37' '  -0.19%  '
Set-TextValue 'D38' '57.27'
Set-TextValue 'E38' '  -0.99%  '
Set-TextValue 'E39' '  -0.11%  '
Set-TextValue 'D40' '0.0₃0723'
Set-TextValue 'E40' '  +0.03%  '
Set-TextValue 'D41' '3.30'
Set-TextValue 'E41' '  -2.42%  '
Set-TextValue 'D42' '33.89'
Set-TextValue 'E42' '  +4.19%  '
Set-TextValue 'D43' '2.72'
Set-TextValue 'E43' '  +0.37%  '
Set-TextValue 'E44' '  +1.28%  '
Set-TextValue 'E45' '  -3.73%  '
Set-TextValue 'D46' '0.0420'
Set-TextValue 'E46' '  +0.66%  '
Set-TextValue 'D47' '3.15'
Set-TextValue 'E47' '  -13.26%  '
Set-TextValue 'D48' '3.24'
Set-TextValue 'E48' '  +2.39%  '
Set-TextValue 'E49' '  -1.13%  '
Set-TextValue 'D50' '2.57'
Set-TextValue 'E50' '  -4.15%  '
Set-TextValue 'D51' '0.999'
